$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 23 (pushes Giovanni..Alex down by one, rows 23-33 -> 24-34)
$ws.Rows.Item(23).Insert()

# Fill in Lauren Macaisa's record in the newly-opened row 23
$ws.Range("A23").Value = "Lauren"
$ws.Range("B23").Value = "Macaisa"
$ws.Range("C23").Value = 5
$ws.Range("F23").Value = "macaisa@broadinstitute.org"
$ws.Range("G23").Value = "Lauren is a Research Associate II focused on organizing and optimizing wet lab protocols for the McCarroll & Macosko BICAN lab efforts. Prior to working at the Broad, Lauren worked at Moffitt Cancer Center developing immunotherapy treatments for breast cancer. She graduated from University of North Florida with her B.S. in Behavioral Neuroscience and is pursuing her M.S. in Biotechnology from Northeastern."
$ws.Range("H23").Value = "Data Generation"

# The row that used to be 23 (Giovanni Marrero) landed on row 24 with a bold
# font carried over from the insert operation
$ws.Range("A24").Font.Bold = $true

# Rows 24-28 are still in the "Data Generation" category as Lauren, so their
# C column (ordinal position within the category) shifts up by one to make
# room for her; rows 29+ belong to the next category and are unaffected
for ($r = 24; $r -le 28; $r++) {
  $ws.Range("C$r").Value = $ws.Range("C$r").Value() + 1
}

# Re-run the autofilter so its reference range covers the new last row
$ws.AutoFilterMode = $false
$ws.Range("A1:H34").AutoFilter()

# Keep the hidden _FilterDatabase defined name in sync with the autofilter range
$fd = $wb.Names.Item("_xlnm._FilterDatabase")
$fd.RefersTo = "=Sheet1!`$A`$1:`$H`$34"

# Match the saved selection state
$ws.Range("A1:H23").Select()
